$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (34 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1849.6061
$ws.Range("I15").Value = 1849.6061
$ws.Range("K15").Value = 5548.8183
$ws.Range("M15").Value = -5379.8183
$ws.Range("H86").Value = 2557.3845
$ws.Range("I86").Value = 2226.2856
$ws.Range("J86").Value = 2943.6667
$ws.Range("K86").Value = 2226.2856
$ws.Range("L86").Value = 2943.6667
$ws.Range("M86").Value = -1103.2856
$ws.Range("N86").Value = -5189.6667
$ws.Range("H89").Value = 2557.3845
$ws.Range("I89").Value = 2226.2856
$ws.Range("J89").Value = 2943.6667
$ws.Range("K89").Value = 11131.428
$ws.Range("L89").Value = 14718.3335
$ws.Range("M89").Value = -5515.428
$ws.Range("N89").Value = -25950.3335
$ws.Range("H101").Value = 829.5833
$ws.Range("I101").Value = 837.125
$ws.Range("K101").Value = 2511.375
$ws.Range("M101").Value = -889.375
$ws.Range("H118").Value = 473.6087
$ws.Range("I118").Value = 494.4762
$ws.Range("K118").Value = 1483.4286
$ws.Range("M118").Value = 173.5714
$ws.Range("H132").Value = 4340.095
$ws.Range("I132").Value = 4340.095
$ws.Range("K132").Value = 13020.285
$ws.Range("M132").Value = -10490.285
$ws.Range("H138").Value = 2940.0889
$ws.Range("J138").Value = 3028.8076
$ws.Range("L138").Value = 9086.4228
$ws.Range("N138").Value = -19366.4228

# --- Sheet: ARM (15 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2596.889
$ws.Range("I88").Value = 1831.75
$ws.Range("K88").Value = 1831.75
$ws.Range("M88").Value = -1425.75
$ws.Range("H91").Value = 2596.889
$ws.Range("I91").Value = 1831.75
$ws.Range("K91").Value = 1831.75
$ws.Range("M91").Value = -427.75
$ws.Range("H122").Value = 1322.8889
$ws.Range("I122").Value = 1175.75
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 3527.25
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -1077.25
$ws.Range("N122").Value = -12400

# --- Sheet: BSM (29 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 592495
$ws.Range("I86").Value = 2002996.8
$ws.Range("J86").Value = 4785.9165
$ws.Range("K86").Value = 2002996.8
$ws.Range("L86").Value = 4785.9165
$ws.Range("M86").Value = -2001873.8
$ws.Range("N86").Value = -7031.9165
$ws.Range("H89").Value = 592495
$ws.Range("I89").Value = 2002996.8
$ws.Range("J89").Value = 4785.9165
$ws.Range("K89").Value = 10014984
$ws.Range("L89").Value = 23929.5825
$ws.Range("M89").Value = -10009368
$ws.Range("N89").Value = -35161.5825
$ws.Range("H94").Value = 1874.3684
$ws.Range("I94").Value = 1520.0714
$ws.Range("K94").Value = 1520.0714
$ws.Range("M94").Value = -1069.0714
$ws.Range("H99").Value = 8574.143
$ws.Range("I99").Value = 8464.538
$ws.Range("K99").Value = 8464.538
$ws.Range("M99").Value = -6966.538
$ws.Range("H105").Value = 4678.7144
$ws.Range("I105").Value = 4446.533
$ws.Range("J105").Value = 5259.1665
$ws.Range("K105").Value = 4446.533
$ws.Range("L105").Value = 5259.1665
$ws.Range("M105").Value = -2699.533
$ws.Range("N105").Value = -8753.166499999999

# --- Sheet: CRP (22 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3979.6155
$ws.Range("I16").Value = 4070.2
$ws.Range("K16").Value = 4070.2
$ws.Range("M16").Value = -3783.2
$ws.Range("H58").Value = 5841.852
$ws.Range("I58").Value = 3694.8333
$ws.Range("J58").Value = 7559.467
$ws.Range("K58").Value = 3694.8333
$ws.Range("L58").Value = 7559.467
$ws.Range("M58").Value = -3491.8333
$ws.Range("N58").Value = -7965.467
$ws.Range("H113").Value = 3979.6155
$ws.Range("I113").Value = 4070.2
$ws.Range("K113").Value = 4070.2
$ws.Range("M113").Value = -1900.2
$ws.Range("H136").Value = 5841.852
$ws.Range("I136").Value = 3694.8333
$ws.Range("J136").Value = 7559.467
$ws.Range("K136").Value = 11084.4999
$ws.Range("L136").Value = 22678.401
$ws.Range("M136").Value = -8534.499899999999
$ws.Range("N136").Value = -27778.401

# --- Sheet: CUL (15 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 1271.75
$ws.Range("J29").Value = 2514.5
$ws.Range("L29").Value = 7543.5
$ws.Range("N29").Value = -8097.5
$ws.Range("H46").Value = 9858106
$ws.Range("I46").Value = 34500400
$ws.Range("K46").Value = 103501200
$ws.Range("M46").Value = -103501109
$ws.Range("H113").Value = 1701.25
$ws.Range("I113").Value = 1730
$ws.Range("J113").Value = 1643.75
$ws.Range("K113").Value = 5190
$ws.Range("L113").Value = 4931.25
$ws.Range("M113").Value = -3020
$ws.Range("N113").Value = -9271.25

# --- Sheet: GSM (23 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 1005
$ws.Range("I18").Value = 1005
$ws.Range("K18").Value = 1005
$ws.Range("M18").Value = -712
$ws.Range("H24").Value = 11666
$ws.Range("J24").Value = 12499
$ws.Range("L24").Value = 12499
$ws.Range("N24").Value = -12845
$ws.Range("H102").Value = 2426.3333
$ws.Range("I102").Value = 2776.8667
$ws.Range("K102").Value = 2776.8667
$ws.Range("M102").Value = -1154.8667
$ws.Range("H113").Value = 68585.03
$ws.Range("I113").Value = 85376.03999999999
$ws.Range("K113").Value = 85376.03999999999
$ws.Range("M113").Value = -83206.03999999999
$ws.Range("H126").Value = 7205.6523
$ws.Range("I126").Value = 7578.9443
$ws.Range("J126").Value = 5861.8
$ws.Range("K126").Value = 22736.8329
$ws.Range("L126").Value = 17585.4
$ws.Range("M126").Value = -20266.8329
$ws.Range("N126").Value = -22525.4

# --- Sheet: LTW (36 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6843.3
$ws.Range("I7").Value = 6127.2
$ws.Range("K7").Value = 6127.2
$ws.Range("M7").Value = -6015.2
$ws.Range("H40").Value = 4571.231
$ws.Range("I40").Value = 4571.231
$ws.Range("K40").Value = 4571.231
$ws.Range("M40").Value = -4435.231
$ws.Range("H68").Value = 4462
$ws.Range("H71").Value = 4462
$ws.Range("H82").Value = 1760.4615
$ws.Range("I82").Value = 1154
$ws.Range("J82").Value = 2139.5
$ws.Range("K82").Value = 1154
$ws.Range("L82").Value = 2139.5
$ws.Range("M82").Value = -793
$ws.Range("N82").Value = -2861.5
$ws.Range("H85").Value = 1760.4615
$ws.Range("I85").Value = 1154
$ws.Range("J85").Value = 2139.5
$ws.Range("K85").Value = 1154
$ws.Range("L85").Value = 2139.5
$ws.Range("M85").Value = 94
$ws.Range("N85").Value = -4635.5
$ws.Range("H93").Value = 5961.25
$ws.Range("I93").Value = 6456.357
$ws.Range("K93").Value = 6456.357
$ws.Range("M93").Value = -5208.357
$ws.Range("H122").Value = 4449.0356
$ws.Range("I122").Value = 3533.6316
$ws.Range("K122").Value = 10600.8948
$ws.Range("M122").Value = -8150.8948
$ws.Range("H126").Value = 6843.3
$ws.Range("I126").Value = 6127.2
$ws.Range("K126").Value = 18381.6
$ws.Range("M126").Value = -15911.6

# --- Sheet: WVR (12 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2517.0908
$ws.Range("I107").Value = 1965.3334
$ws.Range("K107").Value = 5896.0002
$ws.Range("M107").Value = -3976.0002
$ws.Range("H113").Value = 1098.4286
$ws.Range("I113").Value = 1221.3334
$ws.Range("K113").Value = 3664.0002
$ws.Range("M113").Value = -1494.0002
$ws.Range("H122").Value = 2813.8333
$ws.Range("I122").Value = 2145.2727
$ws.Range("K122").Value = 6435.8181
$ws.Range("M122").Value = -3985.8181

Write-Host "Applied 186 cell updates across 8 sheets"
